# Update the description for branch #100003 and remove the now-merged
# #100004 "Add Rest endpoints" row, since its work is folded into #100003.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = "Add features for TruckWays, Trucks and Drivers"
$ws.Rows("6:6").Delete()

$ws.Range("B6").Select()
